$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-3: account holder name / card number
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long digit string that must stay TEXT (not be coerced to a
# number) -- format a scratch cell as Text, copy/paste-special the value
# across (values only) so B3's own style/number-format is untouched, then
# clean the scratch cell back up.
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "2570314725427075"
$ws.Range("Z1").Copy()
$ws.Range("B3").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("Z1").Clear()

$ws.Range("C3").Value = "Mohaupt"

# Row 5: opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 06.04.2025"

# Row 6
$ws.Range("B6").Value = "07.04."
$ws.Range("C6").Value = "08.04."
$ws.Range("D6").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 15763350"
$ws.Range("E6").Value = "83,89-"

# Row 7
$ws.Range("B7").Value = "10.04."
$ws.Range("C7").Value = "11.04."
$ws.Range("D7").Value = "BEITRAG Allianz SE K-14452348"
$ws.Range("E7").Value = "53,67-"

# Row 8
$ws.Range("B8").Value = "12.04."
$ws.Range("C8").Value = "13.04."
$ws.Range("D8").Value = "KARTENZAHLUNG ARAL TANKSTELLE"
$ws.Range("E8").Value = "47,67-"

# Rows 9-11: no longer used as transaction rows -> clear them out
$ws.Range("B9:E11").Value = ""

# Re-shape the (now blank) amount cells' alignment to match the new layout:
#  E9  -> centered + vertical-centered + wrap
#  E10 -> right aligned + vertical-centered + wrap
#  E11 -> right aligned + vertical-centered + wrap
$ws.Range("E9").WrapText = $true
$ws.Range("E9").VerticalAlignment = -4108
$ws.Range("E9").HorizontalAlignment = -4108

$ws.Range("E10:E11").WrapText = $true
$ws.Range("E10:E11").VerticalAlignment = -4108

# Row 12: closing balance
$ws.Range("D12").Value = "KONTOSTAND AM 17.04.2025"
$ws.Range("E12").Value = "185,23-"

# Row 13: next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 27.04.2025"
